$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived for "Espárragos" at Vega Modelo de Temuco.
# It is inserted as the new row 3 (most recent date), pushing the existing
# rows 3:20 down to 4:21.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new record's data.
$ws.Cells.Item(3, 1).Value = 10
$ws.Cells.Item(3, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(3, 3).Value = "La Araucanía"
$ws.Cells.Item(3, 4).Value = 44473
$ws.Cells.Item(3, 5).Value = 9
$ws.Cells.Item(3, 6).Value = 300000000
$ws.Cells.Item(3, 7).Value = "Espárragos"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 200
$ws.Cells.Item(3, 11).Value = 1700
$ws.Cells.Item(3, 12).Value = 1700
$ws.Cells.Item(3, 13).Value = 1700
$ws.Cells.Item(3, 14).Value = "$/kilo"
$ws.Cells.Item(3, 15).Value = "Región del Maule"
$ws.Cells.Item(3, 16).Value = 1700
$ws.Cells.Item(3, 17).Value = 1
$ws.Cells.Item(3, 18).Value = "Hortaliza"

# Match the date-number format used by the other rows' "Fecha" column.
$ws.Cells.Item(3, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
